# Update cryptocurrency price/volume snapshot values (data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.930.96"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "1.793.34"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.09"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5140"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3940"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("E9").Value = "  -7.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.093"
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.95"
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.248"
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.000"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.24"
$ws.Range("E14").Value = "  -4.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.243"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").Value = "1.783.10"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.60"
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001077"
$ws.Range("E18").Value = "  -5.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06525"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.10"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.932"
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").Value = "27.998.39"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.229"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.42"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.40"
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("D28").Value = "1.994.02"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.367"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.76"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1082"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.041"
$ws.Range("E32").Value = "  -5.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.630"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.495"
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07084"
$ws.Range("E35").Value = "  -9.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.905"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02306"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2130"
$ws.Range("E38").Value = "  -4.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.54"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.020"
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6113"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.13"
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("E45").Value = "  -6.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5912"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.717"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.86"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.207"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  -4.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06815"
$ws.Range("E51").Value = "  -2.50%  "
